# "Detection adjusted by clinical status"
# Adds a new parameter row ("rel_detection_clinical") to the "constant"
# sheet, mirroring the existing uniform-distribution parameter rows
# (value / distribution / distri_param1 / distri_param2 in columns B:E).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("constant")
$ws.Activate()

$ws.Range("A27").Value = "rel_detection_clinical"
$ws.Range("B27").Value = 2
$ws.Range("C27").Value = "uniform"
$ws.Range("D27").Value = 1
$ws.Range("E27").Value = 5

# Match the author's on-screen selection/scroll position after the edit.
$ws.Range("G13").Select()
$excel.ActiveWindow.ScrollRow = 9
$excel.ActiveWindow.ScrollColumn = 1
